$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Unit" column between "Address" (B) and "Service" (C) ---
$ws.Columns.Item(3).Insert()

# --- Insert 3 extra "Reading from tariff №N" columns after the existing
#     "Reading" column (which is now column G, after the insert above) ---
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(8).Insert()

# --- Row 1 headers ---
$ws.Range("C1").Value = "Unit"
$ws.Range("G1").Value = "Reading from tariff №1"
$ws.Range("H1").Value = "Reading from tariff №2"
$ws.Range("I1").Value = "Reading from tariff №3"
$ws.Range("J1").Value = "Reading from tariff №4"

# --- Row 2 template placeholders ---
$ws.Range("C2").Value = "{d.meter[i].unitName}"
$ws.Range("H2").Value = "{d.meter[i].value2}"
$ws.Range("I2").Value = "{d.meter[i].value3}"
$ws.Range("J2").Value = "{d.meter[i].value4}"

# --- Row 3 template placeholders ---
$ws.Range("C3").Value = "{d.meter[i + 1].unitName}"
$ws.Range("H3").Value = "{d.meter[i + 1].value2}"
$ws.Range("I3").Value = "{d.meter[i + 1].value3}"
$ws.Range("J3").Value = "{d.meter[i + 1].value4}"
